# Adding PC SunEnergy to the Portfolio Forecast
# Shift all timestamps (column A, rows 2-97) forward by 2 days, and
# update the Actual Production (MW) values (column B) for rows 2-44
# with the new production figures. Rows 45-97 already contain 0 and
# remain unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Actual Production (MW) values for rows 2 through 44 (in order).
$newProduction = @(
    1141,1137,1154,1175,1156,1168,1172,1180,1189,1177,
    1149,1070,1013,1009,976,931,898,829,750,705,
    673,646,605,563,504,429,391,381,340,321,
    310,300,296,323,353,371,412,429,455,453,
    457,464,502
)

for ($row = 2; $row -le 97; $row++) {
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = $cellA.Value2 + 2

    if ($row -le 44) {
        $ws.Cells.Item($row, 2).Value = $newProduction[$row - 2]
    }
}
